$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AmortTemplateGrid-Reason")
$ws.Name = "AmortTemplateGrid"
